# The source data for 2025-01-22 (row 3) had a temporary data glitch where
# "greenpeace.eu" (H) and "wwfeu.bsky.social" (I) follower counts were
# recorded as 0. This upload corrects that by clearing those two cells
# back to empty, matching the rest of the sheet's blank-when-unknown
# convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3:I3").ClearContents()

# Leave the selection where the edit was made, as in the uploaded file.
$ws.Range("H3").Select()
